# Apply the crypto price/volume refresh captured in the commit diff.
# D-column prices (and the Maker/VeChain swap in D36/D37) are numeric-looking
# text (e.g. '1.001', '29.399.45') that Excel's Value setter would silently
# coerce into real numbers, so those cells are briefly switched to the Text
# number format while the value is written, then restored to the default
# "Normal" style so no visible formatting change is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '29.399.45' }
    @{ Cell = 'E2'; Value = '  +0.03%  ' }
    @{ Cell = 'D3'; Value = '1.850.63' }
    @{ Cell = 'E3'; Value = '  +0.13%  ' }
    @{ Cell = 'D4'; Value = '0.9998' }
    @{ Cell = 'E4'; Value = '  +0.02%  ' }
    @{ Cell = 'D5'; Value = '240.45' }
    @{ Cell = 'E5'; Value = '  +0.08%  ' }
    @{ Cell = 'D6'; Value = '0.6287' }
    @{ Cell = 'E6'; Value = '  -0.07%  ' }
    @{ Cell = 'D7'; Value = '1.001' }
    @{ Cell = 'E7'; Value = '  +0.03%  ' }
    @{ Cell = 'D8'; Value = '0.07626' }
    @{ Cell = 'E8'; Value = '  +0.01%  ' }
    @{ Cell = 'D9'; Value = '0.2907' }
    @{ Cell = 'E9'; Value = '  -1.11%  ' }
    @{ Cell = 'D10'; Value = '24.74' }
    @{ Cell = 'E10'; Value = '  +1.07%  ' }
    @{ Cell = 'D11'; Value = '0.07744' }
    @{ Cell = 'E11'; Value = '  -0.02%  ' }
    @{ Cell = 'D12'; Value = '5.038' }
    @{ Cell = 'E12'; Value = '  +0.61%  ' }
    @{ Cell = 'D13'; Value = '0.6791' }
    @{ Cell = 'E13'; Value = '  +0.02%  ' }
    @{ Cell = 'D14'; Value = '0.00001053' }
    @{ Cell = 'E14'; Value = '  -3.31%  ' }
    @{ Cell = 'D15'; Value = '83.23' }
    @{ Cell = 'E15'; Value = '  -0.25%  ' }
    @{ Cell = 'D16'; Value = '6.174' }
    @{ Cell = 'E16'; Value = '  +0.71%  ' }
    @{ Cell = 'D17'; Value = '29.415.11' }
    @{ Cell = 'E17'; Value = '  -0.04%  ' }
    @{ Cell = 'D18'; Value = '228.23' }
    @{ Cell = 'E18'; Value = '  -0.07%  ' }
    @{ Cell = 'E19'; Value = '  -0.66%  ' }
    @{ Cell = 'D21'; Value = '7.504' }
    @{ Cell = 'E21'; Value = '  +0.83%  ' }
    @{ Cell = 'D22'; Value = '1.001' }
    @{ Cell = 'E22'; Value = '  +0.01%  ' }
    @{ Cell = 'D23'; Value = '158.83' }
    @{ Cell = 'E23'; Value = '  +1.03%  ' }
    @{ Cell = 'E24'; Value = '  -0.03%  ' }
    @{ Cell = 'D25'; Value = '8.410' }
    @{ Cell = 'E25'; Value = '  +0.46%  ' }
    @{ Cell = 'D26'; Value = '17.72' }
    @{ Cell = 'E26'; Value = '  +0.39%  ' }
    @{ Cell = 'D27'; Value = '1.394' }
    @{ Cell = 'E27'; Value = '  +7.24%  ' }
    @{ Cell = 'D28'; Value = '1.459' }
    @{ Cell = 'E28'; Value = '  -0.52%  ' }
    @{ Cell = 'D29'; Value = '0.05606' }
    @{ Cell = 'E29'; Value = '  -0.44%  ' }
    @{ Cell = 'D30'; Value = '4.114' }
    @{ Cell = 'E30'; Value = '  +0.08%  ' }
    @{ Cell = 'D31'; Value = '4.069' }
    @{ Cell = 'E31'; Value = '  +0.68%  ' }
    @{ Cell = 'D32'; Value = '1.165' }
    @{ Cell = 'E32'; Value = '  +0.76%  ' }
    @{ Cell = 'E33'; Value = '  -0.74%  ' }
    @{ Cell = 'D34'; Value = '0.6982' }
    @{ Cell = 'E34'; Value = '  -1.55%  ' }
    @{ Cell = 'D35'; Value = '2.582' }
    @{ Cell = 'E35'; Value = '  -0.17%  ' }
    @{ Cell = 'B36'; Value = 'Maker' }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' }
    @{ Cell = 'D36'; Value = '1.234.72' }
    @{ Cell = 'E36'; Value = '  +0.42%  ' }
    @{ Cell = 'B37'; Value = 'VeChain' }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D37'; Value = '0.01808' }
    @{ Cell = 'E37'; Value = '  +0.70%  ' }
    @{ Cell = 'D38'; Value = '2.715' }
    @{ Cell = 'E38'; Value = '  -2.11%  ' }
    @{ Cell = 'D39'; Value = '6.389' }
    @{ Cell = 'E39'; Value = '  -1.48%  ' }
    @{ Cell = 'D40'; Value = '0.9041' }
    @{ Cell = 'E40'; Value = '  -0.37%  ' }
    @{ Cell = 'D41'; Value = '1.001' }
    @{ Cell = 'E41'; Value = '  +0.09%  ' }
    @{ Cell = 'D42'; Value = '101.55' }
    @{ Cell = 'E42'; Value = '  +0.09%  ' }
    @{ Cell = 'D43'; Value = '66.09' }
    @{ Cell = 'E43'; Value = '  +0.09%  ' }
    @{ Cell = 'D44'; Value = '7.220' }
    @{ Cell = 'E44'; Value = '  +0.99%  ' }
    @{ Cell = 'E45'; Value = '  -2.84%  ' }
    @{ Cell = 'D46'; Value = '0.4005' }
    @{ Cell = 'E46'; Value = '  -0.10%  ' }
    @{ Cell = 'D47'; Value = '9.031' }
    @{ Cell = 'E47'; Value = '  +0.24%  ' }
    @{ Cell = 'D48'; Value = '1.680' }
    @{ Cell = 'E48'; Value = '  -0.22%  ' }
    @{ Cell = 'E49'; Value = '  +1.12%  ' }
    @{ Cell = 'D50'; Value = '0.05706' }
    @{ Cell = 'E50'; Value = '  -0.04%  ' }
    @{ Cell = 'D51'; Value = '0.4631' }
    @{ Cell = 'E51'; Value = '  +0.13%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"   # force text so digit-only strings are not parsed as numbers
    $rng.Value = $u.Value
    $rng.Style = "Normal"     # restore the original (unstyled) cell style
}
